$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.189.11"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "2.477.73"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'576.91"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'146.48"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "2.478.33"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "'29.31"
$ws.Range("E14").Value = "  +10.10%  "
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "2.918.66"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "63.237.49"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "2.479.90"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "'7.94"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'11.12"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").Value = "'330.64"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'2.22"
$ws.Range("E23").Value = "  +9.84%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'66.36"
$ws.Range("D26").Value = "'666.45"
$ws.Range("E26").Value = "  +8.88%  "
$ws.Range("D27").Value = "'9.18"
$ws.Range("E27").Value = "  +9.25%  "
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +4.32%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "'4.80"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "'5.53"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.374"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'153.46"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'18.87"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("D42").Value = "'2.74"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").Value = "'1.77"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "0.0₆0300"
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("D46").Value = "'15.13"
$ws.Range("E46").Value = "  +27.46%  "
$ws.Range("D47").Value = "'148.06"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "'20.91"
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("E51").Value = "  +1.15%  "
